$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Ephb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.071327
$ws.Range("H2").Value = 0.213981
$ws.Range("I2").Value = 0.03356605248408491
$ws.Range("J2").Value = 0.03356605248408491
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.095455
$ws.Range("N2").Value = 0.286365
$ws.Range("O2").Value = 0.03201851307709132
$ws.Range("P2").Value = 0.03201851307709131
$ws.Range("Q2").Value = 0.006808518785
$ws.Range("R2").Value = 0.061276669065
$ws.Range("S2").Value = 0.001074735090408006
$ws.Range("T2").Value = 0.001074735090408006

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Ephb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.071327
$ws.Range("H3").Value = 0.213981
$ws.Range("I3").Value = 0.03356605248408491
$ws.Range("J3").Value = 0.03356605248408491
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.220310333333333
$ws.Range("N3").Value = 3.660931
$ws.Range("O3").Value = 0.4093292375039861
$ws.Range("P3").Value = 0.409329237503986
$ws.Range("Q3").Value = 0.08704107514566665
$ws.Range("R3").Value = 0.783369676311
$ws.Range("S3").Value = 0.01373956666932926
$ws.Range("T3").Value = 0.01373956666932925

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Ephb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.071327
$ws.Range("H4").Value = 0.213981
$ws.Range("I4").Value = 0.03356605248408491
$ws.Range("J4").Value = 0.03356605248408491
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.665478666666666
$ws.Range("N4").Value = 4.996435999999999
$ws.Range("O4").Value = 0.5586522494189227
$ws.Range("P4").Value = 0.5586522494189227
$ws.Range("Q4").Value = 0.1187935968573333
$ws.Range("R4").Value = 1.069142371716
$ws.Range("S4").Value = 0.01875175072434765
$ws.Range("T4").Value = 0.01875175072434765

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Ephb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.308223
$ws.Range("H5").Value = 0.9246690000000001
$ws.Range("I5").Value = 0.1450478695977975
$ws.Range("J5").Value = 0.1450478695977975
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.095455
$ws.Range("N5").Value = 0.286365
$ws.Range("O5").Value = 0.03201851307709132
$ws.Range("P5").Value = 0.03201851307709131
$ws.Range("Q5").Value = 0.029421426465
$ws.Range("R5").Value = 0.264792838185
$ws.Range("S5").Value = 0.004644217109521316
$ws.Range("T5").Value = 0.004644217109521315

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Ephb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.308223
$ws.Range("H6").Value = 0.9246690000000001
$ws.Range("I6").Value = 0.1450478695977975
$ws.Range("J6").Value = 0.1450478695977975
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.220310333333333
$ws.Range("N6").Value = 3.660931
$ws.Range("O6").Value = 0.4093292375039861
$ws.Range("P6").Value = 0.409329237503986
$ws.Range("Q6").Value = 0.376127711871
$ws.Range("R6").Value = 3.385149406839
$ws.Range("S6").Value = 0.05937233386404406
$ws.Range("T6").Value = 0.05937233386404405

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Ephb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.308223
$ws.Range("H7").Value = 0.9246690000000001
$ws.Range("I7").Value = 0.1450478695977975
$ws.Range("J7").Value = 0.1450478695977975
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.665478666666666
$ws.Range("N7").Value = 4.996435999999999
$ws.Range("O7").Value = 0.5586522494189227
$ws.Range("P7").Value = 0.5586522494189227
$ws.Range("Q7").Value = 0.513338831076
$ws.Range("R7").Value = 4.620049479684
$ws.Range("S7").Value = 0.08103131862423214
$ws.Range("T7").Value = 0.08103131862423214

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Ephb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.745424333333333
$ws.Range("H8").Value = 5.236273
$ws.Range("I8").Value = 0.8213860779181176
$ws.Range("J8").Value = 0.8213860779181176
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.095455
$ws.Range("N8").Value = 0.286365
$ws.Range("O8").Value = 0.03201851307709132
$ws.Range("P8").Value = 0.03201851307709131
$ws.Range("Q8").Value = 0.1666094797383333
$ws.Range("R8").Value = 1.499485317645
$ws.Range("S8").Value = 0.026299560877162
$ws.Range("T8").Value = 0.02629956087716199

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Ephb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.745424333333333
$ws.Range("H9").Value = 5.236273
$ws.Range("I9").Value = 0.8213860779181176
$ws.Range("J9").Value = 0.8213860779181176
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.220310333333333
$ws.Range("N9").Value = 3.660931
$ws.Range("O9").Value = 0.4093292375039861
$ws.Range("P9").Value = 0.409329237503986
$ws.Range("Q9").Value = 2.129959350018111
$ws.Range("R9").Value = 19.169634150163
$ws.Range("S9").Value = 0.3362173369706128
$ws.Range("T9").Value = 0.3362173369706127

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb3"
$ws.Range("C10").Value = "Ephb6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.745424333333333
$ws.Range("H10").Value = 5.236273
$ws.Range("I10").Value = 0.8213860779181176
$ws.Range("J10").Value = 0.8213860779181176
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.665478666666666
$ws.Range("N10").Value = 4.996435999999999
$ws.Range("O10").Value = 0.5586522494189227
$ws.Range("P10").Value = 0.5586522494189227
$ws.Range("Q10").Value = 2.906966991447555
$ws.Range("R10").Value = 26.16270292302799
$ws.Range("S10").Value = 0.4588691800703429
$ws.Range("T10").Value = 0.4588691800703429
